$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7823839783668518
$ws.Range("B1").Value = 1.150339961051941
$ws.Range("C1").Value = 2.332760334014893
$ws.Range("D1").Value = 3.868865966796875
$ws.Range("E1").Value = 1.888971090316772
